$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.668.74"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.819.83"
$ws.Range("E3").Value = "  +1.35%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.48"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.583"
$ws.Range("E6").Value = "  +4.91%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "34.87"
$ws.Range("E8").Value = "  +7.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.302"
$ws.Range("E9").Value = "  +1.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0702"
$ws.Range("E10").Value = "  +1.15%  "
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.082.78"
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.43"
$ws.Range("E13").Value = "  +2.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.825.20"
$ws.Range("E14").Value = "  +1.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.645"
$ws.Range("E15").Value = "  +1.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.668.38"
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("E17").Value = "  +1.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.26"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0803"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "246.36"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.59"
$ws.Range("E21").Value = "  +3.74%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "173.22"
$ws.Range("E24").Value = "  +5.86%  "
$ws.Range("E25").Value = "  +1.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.50"
$ws.Range("E26").Value = "  +3.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.84"
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("E28").Value = "  +3.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.00"
$ws.Range("E30").Value = "  +2.28%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.86"
$ws.Range("E31").Value = "  +1.74%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0531"
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.24"
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.59"
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.407.38"
$ws.Range("E36").Value = "  -2.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.680"
$ws.Range("E37").Value = "  +1.45%  "
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "84.20"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  +4.93%  "
$ws.Range("E42").Value = "  +2.01%  "
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.76"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("E45").Value = "  +2.68%  "
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.06"
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.983.46"
$ws.Range("E48").Value = "  +1.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.40"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0131"
$ws.Range("E50").Value = "  -1.05%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  +0.12%  "
